# 1/17; separated data processing to dataproc.r
# Append the new day's COVID data row and rename the sheet to match the
# data it now holds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from the generic default to "covdata".
$ws.Name = "covdata"

# Append the next day of data (2021-01-17 == serial 44213): Tests=23, Positive=0.
$ws.Range("A171").Value = 44213
$ws.Range("B171").Value = 23
$ws.Range("C171").Value = 0

# Leave the selection on the newly entered row, as the author would after typing it in.
[void]$ws.Range("E171").Select()
